# Scheduled market-data refresh: updates cached currentAveragePrice /
# currentAveragePriceNQ / currentAveragePriceHQ / LevePriceNQ / LevePriceHQ /
# LeveProfitNQ / LeveProfitHQ figures (columns H-N) for the affected Leve
# rows on each class sheet (ALC, ARM, BSM, CRP, LTW, WVR). Figures are plain
# cached values (no formulas in this workbook), so each touched cell is
# written directly with the refreshed market value. A couple of rows also
# gain/lose a trailing H..N cell because the refreshed quantities/HQ pricing
# pushed a computed profit cell into/out of existence.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 32
$ws.Range("H32").Value = 1397.5714
$ws.Range("I32").Value = 439.8
$ws.Range("K32").Value = 439.8
$ws.Range("M32").Value = -113.8
# Row 51
$ws.Range("H51").Value = 10218
$ws.Range("J51").Value = 11272.5
$ws.Range("L51").Value = 11272.5
$ws.Range("N51").Value = -12240.5
# Row 64
$ws.Range("H64").Value = 6191.5454
$ws.Range("J64").Value = 6191.5454
$ws.Range("L64").Value = 6191.5454
$ws.Range("N64").Value = -6687.5454
# Row 67
$ws.Range("H67").Value = 6191.5454
$ws.Range("J67").Value = 6191.5454
$ws.Range("L67").Value = 6191.5454
$ws.Range("N67").Value = -7907.5454
# Row 70
$ws.Range("H70").Value = 2424.625
$ws.Range("J70").Value = 4999.5
$ws.Range("L70").Value = 14998.5
$ws.Range("N70").Value = -15538.5
# Row 73
$ws.Range("H73").Value = 2424.625
$ws.Range("J73").Value = 4999.5
$ws.Range("L73").Value = 14998.5
$ws.Range("N73").Value = -16870.5
# Row 86
$ws.Range("H86").Value = 147636.14
$ws.Range("I86").Value = 337263
$ws.Range("J86").Value = 5416
$ws.Range("K86").Value = 337263
$ws.Range("L86").Value = 5416
$ws.Range("M86").Value = -336140
$ws.Range("N86").Value = -7662
# Row 89
$ws.Range("H89").Value = 147636.14
$ws.Range("I89").Value = 337263
$ws.Range("J89").Value = 5416
$ws.Range("K89").Value = 1686315
$ws.Range("L89").Value = 27080
$ws.Range("M89").Value = -1680699
$ws.Range("N89").Value = -38312
# Row 132
$ws.Range("H132").Value = 4590.295
$ws.Range("I132").Value = 3986.8545
$ws.Range("K132").Value = 11960.5635
$ws.Range("M132").Value = -9430.5635

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 24
$ws.Range("H24").Value = 32512.6
$ws.Range("J24").Value = 29642
$ws.Range("L24").Value = 29642
$ws.Range("N24").Value = -30390
# Row 63
$ws.Range("H63").Value = 3351.3044
$ws.Range("I63").Value = 3337.1428
$ws.Range("K63").Value = 3337.1428
$ws.Range("M63").Value = -2651.1428
# Row 66
$ws.Range("H66").Value = 3351.3044
$ws.Range("I66").Value = 3337.1428
$ws.Range("K66").Value = 16685.714
$ws.Range("M66").Value = -13253.714
# Row 82
$ws.Range("H82").Value = 53078
$ws.Range("J82").Value = 53078
$ws.Range("L82").Value = 53078
$ws.Range("N82").Value = -53800
# Row 85
$ws.Range("H85").Value = 53078
$ws.Range("J85").Value = 53078
$ws.Range("L85").Value = 53078
$ws.Range("N85").Value = -55574
# Row 88
$ws.Range("H88").Value = 1548.5
$ws.Range("J88").Value = 1927.8
$ws.Range("L88").Value = 1927.8
$ws.Range("N88").Value = -2739.8
# Row 91
$ws.Range("H91").Value = 1548.5
$ws.Range("J91").Value = 1927.8
$ws.Range("L91").Value = 1927.8
$ws.Range("N91").Value = -4735.8
# Row 100
$ws.Range("H100").Value = 32512.6
$ws.Range("J100").Value = 29642
$ws.Range("L100").Value = 29642
$ws.Range("N100").Value = -31806

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 82
$ws.Range("H82").Value = 20873.125
$ws.Range("I82").Value = 20873.125
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 20873.125
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = ""
$ws.Range("N82").Value = -20490.125
# Row 85
$ws.Range("H85").Value = 20873.125
$ws.Range("I85").Value = 20873.125
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 20873.125
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = ""
$ws.Range("N85").Value = -19547.125
# Row 86
$ws.Range("H86").Value = 3900.6667
$ws.Range("J86").Value = 4999.5
$ws.Range("L86").Value = 4999.5
$ws.Range("N86").Value = -7245.5
# Row 89
$ws.Range("H89").Value = 3900.6667
$ws.Range("J89").Value = 4999.5
$ws.Range("L89").Value = 24997.5
$ws.Range("N89").Value = -36229.5
# Row 100
$ws.Range("H100").Value = 45333
$ws.Range("J100").Value = 45333
$ws.Range("L100").Value = 45333
$ws.Range("N100").Value = -47497
# Row 102
$ws.Range("H102").Value = 12631.167
$ws.Range("I102").Value = 12631.167
$ws.Range("K102").Value = 12631.167
$ws.Range("M102").Value = -9386.166999999999

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 62
$ws.Range("H62").Value = 1800
$ws.Range("J62").Value = 1800
$ws.Range("L62").Value = 1800
$ws.Range("N62").Value = -3048
# Row 65
$ws.Range("H65").Value = 1800
$ws.Range("J65").Value = 1800
$ws.Range("L65").Value = 9000
$ws.Range("N65").Value = -15240
# Row 132
$ws.Range("H132").Value = 3092.5
$ws.Range("I132").Value = 3081.8
$ws.Range("K132").Value = 9245.400000000001
$ws.Range("M132").Value = -6715.400000000001

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 601200.6
$ws.Range("I7").Value = 784355.0600000001
$ws.Range("K7").Value = 784355.0600000001
$ws.Range("M7").Value = -784243.0600000001
# Row 18
$ws.Range("H18").Value = 6670
$ws.Range("J18").Value = 20000
$ws.Range("L18").Value = 20000
$ws.Range("N18").Value = -20344
# Row 46
$ws.Range("H46").Value = 3807.111
$ws.Range("I46").Value = 2132.6667
$ws.Range("K46").Value = 2132.6667
$ws.Range("M46").Value = -1944.6667
# Row 68
$ws.Range("H68").Value = 5537.222
$ws.Range("I68").Value = 3996
$ws.Range("J68").Value = 5977.5713
$ws.Range("K68").Value = 3996
$ws.Range("L68").Value = 5977.5713
$ws.Range("M68").Value = -3247
$ws.Range("N68").Value = -7475.5713
# Row 71
$ws.Range("H71").Value = 5537.222
$ws.Range("I71").Value = 3996
$ws.Range("J71").Value = 5977.5713
$ws.Range("K71").Value = 19980
$ws.Range("L71").Value = 29887.8565
$ws.Range("M71").Value = -16236
$ws.Range("N71").Value = -37375.85649999999
# Row 81
$ws.Range("H81").Value = 53989
$ws.Range("J81").Value = 53989
$ws.Range("L81").Value = 53989
$ws.Range("N81").Value = -55985
# Row 82
$ws.Range("H82").Value = 3498.625
$ws.Range("I82").Value = 2199.8
$ws.Range("J82").Value = 4089
$ws.Range("K82").Value = 2199.8
$ws.Range("L82").Value = 4089
$ws.Range("M82").Value = -1838.8
$ws.Range("N82").Value = -4811
# Row 84
$ws.Range("H84").Value = 53989
$ws.Range("J84").Value = 53989
$ws.Range("L84").Value = 161967
$ws.Range("N84").Value = -171951
# Row 85
$ws.Range("H85").Value = 3498.625
$ws.Range("I85").Value = 2199.8
$ws.Range("J85").Value = 4089
$ws.Range("K85").Value = 2199.8
$ws.Range("L85").Value = 4089
$ws.Range("M85").Value = -951.8000000000002
$ws.Range("N85").Value = -6585
# Row 126
$ws.Range("H126").Value = 601200.6
$ws.Range("I126").Value = 784355.0600000001
$ws.Range("K126").Value = 2353065.18
$ws.Range("M126").Value = -2350595.18

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 119
$ws.Range("H119").Value = 49999.75
$ws.Range("J119").Value = 49999.75
$ws.Range("L119").Value = 49999.75
$ws.Range("N119").Value = -59675.75
Write-Output "Applied all market-price updates across ALC/ARM/BSM/CRP/LTW/WVR sheets."
